# Add a new day (06-03) of carjacking data: update the "through" date shown
# in the sheet name / June row label, add the June-7 row's new day-of count
# in column B, and roll the new counts through the rest of the June row and
# the Total row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Name = "Through 2022-06-03"
$ws.Range("A7").Value = "June (through 06-03)"

# Row 7 (June) - new day of data added
$ws.Range("B7").Value = 1
$ws.Range("D7").Value = 6
$ws.Range("E7").Value = 10
$ws.Range("G7").Value = 26
$ws.Range("H7").Value = 9
$ws.Range("I7").Value = 9

# Row 8 (Total) - updated totals
$ws.Range("B8").Value = 109
$ws.Range("D8").Value = 322
$ws.Range("E8").Value = 305
$ws.Range("G8").Value = 384
$ws.Range("H8").Value = 640
$ws.Range("I8").Value = 673
